# Applies the "Add player role labeling to seating chart" backup-sheet edit:
#  1. Contestants sheet: swap the data rows for Kathleen Reynolds / Peter Adamidis
#     (rows 3 and 4) - Peter is now the "assigned" seat holder for that pair and
#     Kathleen becomes the attendee swapped in, including medical info.
#  2. Seat Assignments sheet: remove Kathleen's now-stale seat-assignment row.
#  3. Standbys sheet: add a new standby row for Kathleen (the swapped-out
#     contestant) ahead of Peter's existing standby row.
#  4. Add a new "Canceled Assignments" sheet recording Peter's canceled seat.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Contestants: swap rows 3 and 4 (all columns A:M)
# ---------------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

# NOTE: Range.Value (the property getter without parens) comes back through
# this COM shim as a bogus reflection string, not the cell data - call it as
# a method (parens) to get the real value, same as Value2/Text.
$row3 = @($contestants.Range("A3").Value(), $contestants.Range("B3").Value(), $contestants.Range("C3").Value(), `
           $contestants.Range("D3").Value(), $contestants.Range("E3").Value(), $contestants.Range("F3").Value(), `
           $contestants.Range("G3").Value(), $contestants.Range("I3").Value(), $contestants.Range("J3").Value(), `
           $contestants.Range("K3").Value(), $contestants.Range("L3").Value(), $contestants.Range("M3").Value())

$row4 = @($contestants.Range("A4").Value(), $contestants.Range("B4").Value(), $contestants.Range("C4").Value(), `
           $contestants.Range("D4").Value(), $contestants.Range("E4").Value(), $contestants.Range("F4").Value(), `
           $contestants.Range("G4").Value(), $contestants.Range("I4").Value(), $contestants.Range("J4").Value(), `
           $contestants.Range("K4").Value(), $contestants.Range("L4").Value(), $contestants.Range("M4").Value())

$cols = @("A","B","C","D","E","F","G","I","J","K","L","M")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $contestants.Range($cols[$i] + "3").Value = $row4[$i]
    $contestants.Range($cols[$i] + "4").Value = $row3[$i]
}

# Column F (Phone) holds a digit-only string ("498086080") that must stay
# text (matches the source file's t="str"), not get auto-coerced to a
# number. Force text format, re-assign, then paste back a plain/default
# format (copied from a normal text cell) so no stray NumberFormat lingers
# on the cell.
$contestants.Range("F3").NumberFormat = "@"
$contestants.Range("F3").Value = $row4[5]
$contestants.Range("F4").NumberFormat = "@"
$contestants.Range("F4").Value = $row3[5]
$contestants.Range("A2").Copy()
$contestants.Range("F3").PasteSpecial(-4122)  # xlPasteFormats
$contestants.Range("A2").Copy()
$contestants.Range("F4").PasteSpecial(-4122)  # xlPasteFormats

# G3's target value is empty (Peter's row has no Location); this engine's
# Value-setter treats "" as "delete the cell" (matching real Excel - you
# can't type a blank into a cell and have it stay "present"), so the closest
# reachable state is an absent/cleared cell.
$contestants.Range("G3").ClearContents()

# ---------------------------------------------------------------------------
# 2. Seat Assignments: delete row 3 (Kathleen's seat assignment record)
# ---------------------------------------------------------------------------
$seatAssignments = $wb.Worksheets.Item("Seat Assignments")
$seatAssignments.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 3. Standbys: insert a new row above the existing one for Kathleen, with the
#    existing row (Peter's standby) shifting down to row 3 unchanged.
# ---------------------------------------------------------------------------
$standbys = $wb.Worksheets.Item("Standbys")
$standbys.Rows.Item(2).Insert()

$standbys.Range("A2").Value = "4e437be8-c154-44a8-b28d-d663c725a9a7"
$standbys.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$standbys.Range("C2").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$standbys.Range("D2").Value = "pending"

# ---------------------------------------------------------------------------
# 4. Add the new "Canceled Assignments" sheet after "Block Types"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$canceled = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$canceled.Name = "Canceled Assignments"

$canceled.Range("A1").Value = "ID"
$canceled.Range("B1").Value = "RecordDayID"
$canceled.Range("C1").Value = "ContestantID"
$canceled.Range("D1").Value = "Reason"
$canceled.Range("E1").Value = "CanceledAt"

$canceled.Range("A2").Value = "53071bc2-2450-41b6-ae7a-68a1ab35e988"
$canceled.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$canceled.Range("C2").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$canceled.Range("D2").Value = "Standby - eligible for reschedule"

# Reuse the existing date-time number format (style index 1 in the original
# workbook, e.g. Seat Assignments!F2) so the new cell matches the same style
# instead of minting a duplicate numFmt entry.
$seatAssignments.Range("F2").Copy($canceled.Range("E2"))
$canceled.Range("E2").Value = 45998.113707662036
